$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.082.94"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.832.71"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'243.37"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D6").Value = "'0.6284"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("E8").Value = "  -1.22%  "

$ws.Range("D9").Value = "'0.2921"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "'23.26"
$ws.Range("E10").Value = "  +3.27%  "

$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "1.834.09"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'5.008"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").Value = "'0.6670"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "'82.74"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "'0.000009378"
$ws.Range("E16").Value = "  -8.63%  "

$ws.Range("D17").Value = "'5.980"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "29.093.17"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").Value = "2.080.60"
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").Value = "'12.58"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("D21").Value = "'223.23"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").Value = "'7.100"
$ws.Range("E23").Value = "  -1.07%  "

$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "'159.84"
$ws.Range("E25").Value = "  +1.07%  "

$ws.Range("D26").Value = "'0.1394"
$ws.Range("E26").Value = "  +1.53%  "

$ws.Range("D27").Value = "'8.490"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'17.89"
$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").Value = "'1.497"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").Value = "'0.05670"
$ws.Range("E30").Value = "  +8.60%  "

$ws.Range("D31").Value = "'4.152"
$ws.Range("E31").Value = "  +1.33%  "

$ws.Range("D32").Value = "'4.084"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("D34").Value = "'1.842"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").Value = "'0.7418"

$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").Value = "'2.669"
$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").Value = "1.221.84"
$ws.Range("E39").Value = "  -1.03%  "

$ws.Range("D40").Value = "'0.01778"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("D41").Value = "'6.532"
$ws.Range("E41").Value = "  +2.79%  "

$ws.Range("D42").Value = "'0.8931"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").Value = "'101.96"
$ws.Range("E44").Value = "  +0.29%  "

$ws.Range("D45").Value = "1.981.46"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'65.88"
$ws.Range("E46").Value = "  +2.94%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000124"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").Value = "'0.5087"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").Value = "'0.4075"
$ws.Range("E49").Value = "  +1.00%  "

$ws.Range("D50").Value = "'0.07442"
$ws.Range("E50").Value = "  +7.13%  "

$ws.Range("D51").Value = "'9.011"
$ws.Range("E51").Value = "  +1.16%  "
